# Switch the active/selected sheet from "ProductLoanOutput" to "ProductLoanInput",
# update the repayment strategy cell value on ProductLoanInput, and give it
# a left/top aligned style (new cell style gets auto-created by the engine).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")

# Make ProductLoanInput the active sheet (was ProductLoanOutput before).
$ws1.Activate()

# Update the "repaymentstrategy" value cell (B17) to the new scenario text.
$ws1.Range("B17").Value = "Penalties, Fees, Interest, Principal order"

# Give the updated cell a left/top alignment (creates a new cell style).
$ws1.Range("B17").HorizontalAlignment = -4131
$ws1.Range("B17").VerticalAlignment = -4160

# Move the selection/active cell to B17 on the now-active sheet.
$ws1.Range("B17").Select()
